# Apply updated pl_mw.xlsx values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each row: row number, then new values for columns B, C, D, F, G, I, L
$data = @(
    @(2,1.714699423880404,0.250909639502737,0.0769704405247964,1.994199913594798,0.002492051564827234,1.469259451262516,0.3037662645572681),
    @(3,1.598190042741692,0.2181842122574267,0.07715322995481344,1.937835160850398,0.002497890905039584,1.446692224825796,0.2928482977020366),
    @(4,1.527687279535769,0.198095409460592,0.07729729758906601,1.904505472940002,0.002501662239036914,1.433679701155455,0.2863377924453516),
    @(5,1.499215713705098,0.189909547437253,0.07736393689206267,1.891241731045469,0.002503246017194831,1.428587366770017,0.2837330116229282),
    @(6,1.494503632703697,0.1885503004507143,0.07737547911300169,1.889058435354471,0.00250351184201403,1.427754441656774,0.283303399306817),
    @(7,1.527302255447466,0.1979850109984227,0.07729816429796088,1.904325308512725,0.002501683407992819,1.433610174774245,0.2863024682353483),
    @(8,1.674311480527138,0.2396244914069428,0.07702681424930446,1.974498095576536,0.002494026476373,1.461302026612586,0.299961488222948),
    @(9,1.970879677782477,0.3213530683602812,0.0767506654778245,2.122400591792058,0.002480479005819775,1.522387282233552,0.3282936770582836),
    @(10,2.193950921345674,0.3815001085517338,0.0767085342868441,2.237560365112188,0.002471409615084408,1.571526071065549,0.3500747776921287),
    @(11,2.296588507953288,0.4088971872736238,0.07672529932866468,2.291408750250667,0.002467473331525527,1.594833016531382,0.3601982950723084),
    @(12,2.335623807493221,0.4192778376123556,0.07673689872933664,2.312013652180383,0.002466009824880642,1.603797998822273,0.3640631207529736),
    @(13,2.327209328130209,0.4170419001647474,0.07673416571185498,2.307566457819121,0.00246632381554787,1.601861008133056,0.3632293654137442),
    @(14,2.299796575246489,0.4097510852453183,0.07672614789994014,2.293099625098165,0.002467352386347366,1.59556777053163,0.3605156280295319),
    @(15,2.283027489638357,0.4052860525773667,0.07672192307979486,2.284266205458408,0.002467985937156297,1.591731162844582,0.3588574675154348),
    @(16,2.187266824357437,0.3797104360768913,0.07670816883833709,2.234070957072078,0.00247167065845815,1.57002226752725,0.3494175398053159),
    @(17,2.128819375527996,0.3640303868009482,0.07670899334465275,2.203654638740716,0.002473979518038638,1.556950214546092,0.3436817937426042),
    @(18,2.095311156999685,0.3550149656264239,0.07671284448196758,2.186297340159712,0.00247532535238376,1.54952113682279,0.3404029786625671),
    @(19,2.083984556843973,0.3519630461830161,0.07671472584784311,2.180443943125681,0.002475784097587946,1.547021109427547,0.3392962936689941),
    @(20,2.135029891734177,0.365699204517,0.07670855547219446,2.206878265901224,0.002473731891004978,1.558332466705778,0.3442902771858058),
    @(21,2.307843773842421,0.4118924057155482,0.07672835974643988,2.297343056996795,0.002467049536530151,1.597412454170353,0.3613118672068083),
    @(22,2.421771619525884,0.4421175267117405,0.07677194846279178,2.357714215533377,0.002462839990518547,1.623765506431212,0.3726188784455928),
    @(23,2.360875578221169,0.4259823181678826,0.07674585231223574,2.325377705206847,0.002465072322774041,1.609625407059553,0.3665673128990932),
    @(24,2.132221824698263,0.3649447344315604,0.07670874292134755,2.205420461124902,0.002473843785603006,1.557707282103465,0.3440151234696884),
    @(25,1.889749850323426,0.29922969188911,0.07679748333037395,2.08126436099991,0.00248398794590211,1.505123201536918,0.3204611711356335)
)

$columns = @("B", "C", "D", "F", "G", "I", "L")

foreach ($row in $data) {
    $rowNum = $row[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $value = $row[$i + 1]
        $ws.Range("$col$rowNum").Value = $value
    }
}
